$d = $word.ActiveDocument

# 1) Simple text replacement: "equipamentos" -> "ferramentas"
$d.Content.Find.Execute(
    "O jogo deve permitir ao jogador visualizar seus equipamentos.  ",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "O jogo deve permitir ao jogador visualizar suas ferramentas.  ",
    2) | Out-Null

# 2) Append new paragraphs (Requisito Funcional 11-15) at the end of the
#    document body, mirroring the existing "Requisito Funcional N: <texto>"
#    pattern (blank paragraph, then a paragraph with a bold label run
#    followed by a normal-text run).

function Add-EmptyParagraph {
    $end = $d.Content
    $end.Collapse(0)
    $end.InsertParagraphAfter() | Out-Null
}

function Add-RequisitoParagraph([string]$label, [string]$body) {
    $end = $d.Content
    $end.Collapse(0)
    $end.InsertParagraphAfter() | Out-Null

    $p = $d.Paragraphs.Item($d.Paragraphs.Count)
    $r = $p.Range
    $start = $r.Start
    $full = $label + $body
    $r.InsertAfter($full) | Out-Null

    $boldEnd = $start + $label.Length
    $boldRange = $d.Range($start, $boldEnd)
    $boldRange.Font.Bold = 1
}

Add-EmptyParagraph
Add-RequisitoParagraph "Requisito Funcional 11: " "O jogo deve permitir ao jogador visualizar a distância percorrida no jogo em número de passos.  "

Add-EmptyParagraph
Add-RequisitoParagraph "Requisito Funcional 12: " "O jogo deve permitir ao jogador visualizar quanto dinheiro ele possui no momento.  "

Add-EmptyParagraph
Add-RequisitoParagraph "Requisito Funcional 13: " "O jogo deve permitir ao jogador contratar funcionários para sua empresa.   "

Add-EmptyParagraph
Add-RequisitoParagraph "Requisito Funcional 14: " "O jogo deve permitir que o jogador veja a descrição dos funcionários antes de serem contratados."

Add-EmptyParagraph

# Final paragraph only contains the bold "Requisito Funcional 15: " label
# (no body text yet), so the paragraph mark itself also ends up bold.
$end = $d.Content
$end.Collapse(0)
$end.InsertParagraphAfter() | Out-Null
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$r = $p.Range
$label = "Requisito Funcional 15: "
$r.InsertAfter($label) | Out-Null
# Bold the whole paragraph (including its end-of-paragraph mark), since the
# entire paragraph content is the bold label.
$p.Range.Font.Bold = 1

Write-Output "done"
